# "SQL for inserting photos"
#
# Adds a new worksheet "20201028" (continuing the date-named sheet series)
# after the last existing sheet, containing the ID / Shop ID / SQL table
# for restaurant da04f5c9-ffb0-11ea-ba65-065a10bcba76, photo ids 160-169,
# with the same CONCAT(...) INSERT-statement-builder formula used on every
# other sheet in the workbook. The new sheet becomes the active tab, and
# the previously-active sheet's selection reverts to the "ready to paste
# the next batch" A1:C2 state used by the other finished sheets.

$wb = $excel.ActiveWorkbook

# The sheet that is currently last / currently active ("20201026").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Insert the new sheet right after it, so it lands at the end of the tab
# strip and becomes the active sheet.
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "20201028"

# Once a new sheet is active, the old sheet's lingering selection is reset
# to A1:C2 (matching every other already-completed sheet in the book).
$lastSheet.Range("A1:C2").Select()

# --- Header row ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Shop ID"
$ws.Range("C1").Value = "SQL"
$ws.Range("A1:C1").Font.Color = 0

# --- Data rows: photo ids 160-169 for restaurant da04f5c9-... ---
$uuid = "da04f5c9-ffb0-11ea-ba65-065a10bcba76"
$startId = 160
for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $startId + $i
    $ws.Range("B$row").Value = $uuid
}

# Shop ID column is always rendered in the explicit-black font style.
$ws.Range("B2:B11").Font.Color = 0
# ID column alternates that same style on the even data rows (2,4,6,8,10).
$ws.Range("A2").Font.Color = 0
$ws.Range("A4").Font.Color = 0
$ws.Range("A6").Font.Color = 0
$ws.Range("A8").Font.Color = 0
$ws.Range("A10").Font.Color = 0

# --- SQL column formula ---
# First row is entered on its own; the rest are entered as one range so the
# relative formula fills down (same pattern as every other sheet's C column).
$ws.Range("C2").Formula = "=_xlfn.CONCAT(""INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('"", B2, ""'), LPAD("", A2, "", 7, '0'), 'dish'"", "");"")"
$ws.Range("C3:C11").Formula = "=_xlfn.CONCAT(""INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('"", B3, ""'), LPAD("", A3, "", 7, '0'), 'dish'"", "");"")"

# Leave the selection where the author left it on the new sheet.
$ws.Range("I6").Select()
